# Exercise 3: add a new "empty list" slide after the existing slide and
# file it into its own section, as per commit "doc: add todos for exercise 3".

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# --- 1. New blank slide (uses the same "Leer" / blank layout as slide 1) ---
$s2 = $p.Slides.Add(2, 12)   # 12 = ppLayoutBlank

# --- 2. Title textbox "Meine kleine Bibliothek" (identical to slide 1's) ---
$s1.Shapes.Item(1).Copy()
$s2.Shapes.Paste()

# --- 3. Rounded-rectangle "empty state" banner, styled like the book rows
#        on slide 1, but repositioned and re-worded ---
$s1.Shapes.Item(2).Copy()
$s2.Shapes.Paste()
$rect = $s2.Shapes.Item($s2.Shapes.Count)
$rect.Name = "Abgerundetes Rechteck 11"
$rect.Left = 46.34480484960632
$rect.Top = 87.5772440944882
$rect.TextFrame.TextRange.Text = "Es sind noch keine Bücher vorhanden"
$rect.TextFrame.TextRange.ParagraphFormat.Alignment = 2   # ppAlignCenter

# --- 4. Copyright footer textbox (identical to slide 1's) ---
$s1.Shapes.Item($s1.Shapes.Count).Copy()
$s2.Shapes.Paste()

# --- 5. File the new slide into its own section ---
$secIdx = $p.SectionProperties.AddSection($p.SectionProperties.Count + 1, "Liste ohne Elemente")
$s2.sectionIndex = $secIdx
$s2.MoveToSectionStart($secIdx)
